$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume snapshot figures (Price = column D, Volume(1h) = column E).
# For Price cells whose new text looks like a plain decimal number (e.g. "213.04"),
# force the cell to Text format first so Excel keeps the literal string instead of
# silently converting it to a floating point number (which would also introduce
# binary floating point rounding noise, e.g. 213.04 -> 213.03999999999999).
$ws.Range("D2").Value = "27.963.88"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.640.10"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.04"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("D12").Value = "1.872.61"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "1.630.01"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("E14").Value = "  +3.82%  "
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.93"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "27.966.23"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.41"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.70"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.32"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").Value = "1.407.88"
$ws.Range("E34").Value = "  -5.11%  "
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.558"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.905"
$ws.Range("E40").Value = "  -5.57%  "
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.39"
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "1.781.40"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.03"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.61"
$ws.Range("E51").Value = "  -1.74%  "
